$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap the "Valor Mora" values between period 2201 (row 16) and period 2212 (row 27)
$ws.Range("F16").Value = 25333
$ws.Range("F27").Value = 32000
